$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-03-05 Tuesday" "2024-03-06 Wednesday"

Replace-Text "76×65=" "57×31="
Replace-Text "34×17=" "55×80="
Replace-Text "41×12=" "16×49="
Replace-Text "95×96=" "58×32="
Replace-Text "48×40=" "48×76="
Replace-Text "85×80=" "73×18="
Replace-Text "60×22=" "11×31="
Replace-Text "21×40=" "83×40="
Replace-Text "44×36=" "63×73="
Replace-Text "82×80=" "69×40="
Replace-Text "74×43=" "66×73="
Replace-Text "52×14=" "18×40="
Replace-Text "20×18=" "77×45="
Replace-Text "20×15=" "63×42="
Replace-Text "81×87=" "40×46="
Replace-Text "18×74=" "19×84="
Replace-Text "59×72=" "39×64="
Replace-Text "39×35=" "46×22="
Replace-Text "70×80=" "77×14="
Replace-Text "81×28=" "24×98="
Replace-Text "84×86=" "83×75="
Replace-Text "75×57=" "70×69="
Replace-Text "60×88=" "18×87="
Replace-Text "74×50=" "62×34="
Replace-Text "53×68=" "65×51="
